$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> "2016-08-17 20:29:28"
foreach ($r in 4..7) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-08-17 20:29:28"
}

# de-de sheet: rows 4-7, Priority (E) -> "ht", Latest Handoff Datetime (H) -> "2016-08-17 20:29:33"
foreach ($r in 4..7) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-08-17 20:29:33"
}
